$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "20.580.88"

$ws.Range("E2").Value = "  +1.21%  "

$ws.Range("D3").Value = "1.475.87"

$ws.Range("E3").Value = "  +0.82%  "

$ws.Range("E4").Value = "  -0.39%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9500"
$ws.Range("D5").ClearFormats()

$ws.Range("E5").Value = "  +5.82%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "279.38"
$ws.Range("D6").ClearFormats()

$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3651"
$ws.Range("D7").ClearFormats()

$ws.Range("E7").Value = "  -1.59%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3066"
$ws.Range("D8").ClearFormats()

$ws.Range("E8").Value = "  -2.90%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.97"
$ws.Range("D9").ClearFormats()

$ws.Range("E9").Value = "  +1.64%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.059"
$ws.Range("D10").ClearFormats()

$ws.Range("E10").Value = "  +2.44%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06669"
$ws.Range("D11").ClearFormats()

$ws.Range("E11").Value = "  +2.05%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.005"
$ws.Range("D12").ClearFormats()

$ws.Range("E12").Value = "  -0.27%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.521"
$ws.Range("D13").ClearFormats()

$ws.Range("E13").Value = "  +0.96%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.07"
$ws.Range("D14").ClearFormats()

$ws.Range("E14").Value = "  +2.47%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.221"
$ws.Range("D15").ClearFormats()

$ws.Range("E15").Value = "  +1.49%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9511"
$ws.Range("D16").ClearFormats()

$ws.Range("E16").Value = "  +5.06%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001032"
$ws.Range("D17").ClearFormats()

$ws.Range("E17").Value = "  +0.88%  "

$ws.Range("D18").Value = "1.473.59"

$ws.Range("E18").Value = "  +0.38%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.05945"
$ws.Range("D19").ClearFormats()

$ws.Range("E19").Value = "  +5.34%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.80"
$ws.Range("D20").ClearFormats()

$ws.Range("E20").Value = "  +3.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.497"
$ws.Range("D21").ClearFormats()

$ws.Range("E21").Value = "  -1.61%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.46"
$ws.Range("D22").ClearFormats()

$ws.Range("E22").Value = "  +0.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.09"
$ws.Range("D23").ClearFormats()

$ws.Range("E23").Value = "  -0.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.261"
$ws.Range("D24").ClearFormats()

$ws.Range("E24").Value = "  -1.20%  "

$ws.Range("D25").Value = "20.614.09"

$ws.Range("E25").Value = "  +0.72%  "

$ws.Range("E26").Value = "  +5.16%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.122"
$ws.Range("D27").ClearFormats()

$ws.Range("E27").Value = "  -4.48%  "

$ws.Range("E28").Value = "  +0.97%  "

$ws.Range("D29").Value = "1.635.60"

$ws.Range("E29").Value = "  +0.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "113.77"
$ws.Range("D30").ClearFormats()

$ws.Range("E30").Value = "  +2.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.960"
$ws.Range("D31").ClearFormats()

$ws.Range("E31").Value = "  +1.26%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.022"
$ws.Range("D32").ClearFormats()

$ws.Range("E32").Value = "  +1.58%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8099"
$ws.Range("D33").ClearFormats()

$ws.Range("E33").Value = "  -0.57%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07973"
$ws.Range("D34").ClearFormats()

$ws.Range("E34").Value = "  +2.84%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.518"
$ws.Range("D35").ClearFormats()

$ws.Range("E35").Value = "  +4.98%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.222"
$ws.Range("D36").ClearFormats()

$ws.Range("E36").Value = "  +6.17%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05845"
$ws.Range("D37").ClearFormats()

$ws.Range("E37").Value = "  -2.58%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.733"
$ws.Range("D38").ClearFormats()

$ws.Range("E38").Value = "  -1.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02057"
$ws.Range("D39").ClearFormats()

$ws.Range("E39").Value = "  +1.58%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.39"
$ws.Range("D40").ClearFormats()

$ws.Range("E40").Value = "  +0.80%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9528"
$ws.Range("D41").ClearFormats()

$ws.Range("E41").Value = "  +3.63%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1880"
$ws.Range("D42").ClearFormats()

$ws.Range("E42").Value = "  +2.34%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.411"
$ws.Range("D43").ClearFormats()

$ws.Range("E43").Value = "  +7.38%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5313"
$ws.Range("D44").ClearFormats()

$ws.Range("E44").Value = "  +0.72%  "

$ws.Range("E45").Value = "  -0.62%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.26"
$ws.Range("D46").ClearFormats()

$ws.Range("E46").Value = "  +1.60%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "118.16"
$ws.Range("D47").ClearFormats()

$ws.Range("E47").Value = "  -2.40%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5196"
$ws.Range("D48").ClearFormats()

$ws.Range("E48").Value = "  +0.33%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.818"
$ws.Range("D49").ClearFormats()

$ws.Range("E49").Value = "  +0.81%  "

$ws.Range("E50").Value = "  +1.72%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9835"
$ws.Range("D51").ClearFormats()

$ws.Range("E51").Value = "  -1.26%  "
